$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of (old type, old value) -> new type code. New value is always "Cheer/Drill".
$map = @{
    "club-sports|Team Tennis-Coed"  = "sports_club_coed"
    "uil-sports|Team Tennis-Coed"   = "sports_uil_coed"
    "uil-sports|Team Tennis-Boys"   = "sports_uil_boys"
    "uil-sports|Team Tennis-Girls"  = "sports_uil_girls"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 29 }

for ($r = 2; $r -le $lastRow; $r++) {
    $typeCell = $ws.Cells.Item($r, 2)
    $valueCell = $ws.Cells.Item($r, 3)

    $oldType = $typeCell.Value2
    $oldValue = $valueCell.Value2

    if ($null -eq $oldType -or $oldType -eq "") { continue }

    $key = "$oldType|$oldValue"
    if ($map.ContainsKey($key)) {
        $typeCell.Value = $map[$key]
        $valueCell.Value = "Cheer/Drill"
    }
}
